$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 16; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    $hVal = $ws.Cells.Item($r, 8).Value2

    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
    $ws.Cells.Item($r, 6).Value = $hVal
    $ws.Cells.Item($r, 8).Value = $fVal
}
